$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.87'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.72%'
$ws.Range("E2").NumberFormat = "General"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.97'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.88%'
$ws.Range("E3").NumberFormat = "General"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.55%'
$ws.Range("E4").NumberFormat = "General"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05616'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.47%'
$ws.Range("E5").NumberFormat = "General"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.497'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.14%'
$ws.Range("E6").NumberFormat = "General"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8181'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.08%'
$ws.Range("E7").NumberFormat = "General"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8322'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.29%'
$ws.Range("E8").NumberFormat = "General"

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'One'
$ws.Range("B9").NumberFormat = "General"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("C9").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0006011'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.72%'
$ws.Range("E9").NumberFormat = "General"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("B10").NumberFormat = "General"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("C10").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1328'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.92%'
$ws.Range("E10").NumberFormat = "General"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("B11").NumberFormat = "General"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("C11").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06946'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.03%'
$ws.Range("E11").NumberFormat = "General"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("B12").NumberFormat = "General"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("C12").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02899'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.67%'
$ws.Range("E12").NumberFormat = "General"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("B13").NumberFormat = "General"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("C13").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09379'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.18%'
$ws.Range("E13").NumberFormat = "General"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("B14").NumberFormat = "General"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C14").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001514'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.07%'
$ws.Range("E14").NumberFormat = "General"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006165'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.69%'
$ws.Range("E15").NumberFormat = "General"

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.16%'
$ws.Range("E16").NumberFormat = "General"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.022'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.17%'
$ws.Range("E17").NumberFormat = "General"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.301'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '8.66%'
$ws.Range("E18").NumberFormat = "General"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03070'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-5.07%'
$ws.Range("E20").NumberFormat = "General"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1292'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.09%'
$ws.Range("E21").NumberFormat = "General"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.742'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.34%'
$ws.Range("E22").NumberFormat = "General"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04586'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.26%'
$ws.Range("E23").NumberFormat = "General"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.42%'
$ws.Range("E24").NumberFormat = "General"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.61%'
$ws.Range("E25").NumberFormat = "General"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004489'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.56%'
$ws.Range("E26").NumberFormat = "General"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009801'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.14%'
$ws.Range("E27").NumberFormat = "General"

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.73%'
$ws.Range("E28").NumberFormat = "General"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03638'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.54%'
$ws.Range("E40").NumberFormat = "General"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006100'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '80.18%'
$ws.Range("E41").NumberFormat = "General"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-22.98%'
$ws.Range("E42").NumberFormat = "General"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002599'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.95%'
$ws.Range("E43").NumberFormat = "General"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008101'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.62%'
$ws.Range("E44").NumberFormat = "General"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005295'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.16%'
$ws.Range("E45").NumberFormat = "General"

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.05%'
$ws.Range("E46").NumberFormat = "General"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-18.31%'
$ws.Range("E47").NumberFormat = "General"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002594'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '21.98%'
$ws.Range("E48").NumberFormat = "General"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.05%'
$ws.Range("E49").NumberFormat = "General"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.05%'
$ws.Range("E50").NumberFormat = "General"

Write-Host "Applied all cell updates"